# Generate Report for Handback
# Update status + handback timestamps for the zh-cn and de-de worksheets,
# and reflect the new status on the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: 9ca4079f... file has now been handed back, and both files'
# Latest Handback DateTime advance to the new handback timestamp.
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("H2").Value = "2016-03-18 17:28:55"
$wsZhCn.Range("H3").Value = "2016-03-18 17:28:55"

# de-de sheet: same update, with its own handback timestamp.
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("H2").Value = "2016-03-18 17:29:09"
$wsDeDe.Range("H3").Value = "2016-03-18 17:29:09"

# Overview sheet reflects the 9ca4079f... row now being handed back for
# both locales.
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"
